$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.250.72'
$ws.Range('E2').Value = '  +2.30%  '

# Row 3
$ws.Range('D3').Value = '2.099.19'
$ws.Range('E3').Value = '  +4.37%  '

# Row 5
$ws.Range('D5').Value = '251.00'
$ws.Range('E5').Value = '  +2.36%  '

# Row 6
$ws.Range('E6').Value = '  +0.52%  '

# Row 7
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('D8').Value = '54.28'
$ws.Range('E8').Value = '  +20.64%  '

# Row 9
$ws.Range('D9').Value = '61.72'
$ws.Range('E9').Value = '  +2.05%  '

# Row 10
$ws.Range('E10').Value = '  +1.37%  '

# Row 11
$ws.Range('E11').Value = '  +4.24%  '

# Row 12
$ws.Range('E12').Value = '  +7.35%  '

# Row 13
$ws.Range('D13').Value = '15.36'
$ws.Range('E13').Value = '  +5.67%  '

# Row 14
$ws.Range('D14').Value = '2.403.92'
$ws.Range('E14').Value = '  +4.20%  '

# Row 15
$ws.Range('D15').Value = '0.840'
$ws.Range('E15').Value = '  +3.92%  '

# Row 16
$ws.Range('D16').Value = '2.100.76'
$ws.Range('E16').Value = '  +4.45%  '

# Row 17
$ws.Range('D17').Value = '5.16'
$ws.Range('E17').Value = '  +5.26%  '

# Row 18
$ws.Range('D18').Value = '37.247.08'
$ws.Range('E18').Value = '  +2.68%  '

# Row 19
$ws.Range('D19').Value = '72.60'
$ws.Range('E19').Value = '  +1.82%  '

# Row 20
$ws.Range('D20').Value = '14.60'
$ws.Range('E20').Value = '  +13.61%  '

# Row 21
$ws.Range('E21').Value = '  +2.73%  '

# Row 22
$ws.Range('D22').Value = '241.66'
$ws.Range('E22').Value = '  +2.26%  '

# Row 23
$ws.Range('E23').Value = '  +7.46%  '

# Row 24
$ws.Range('E24').Value = '  +0.10%  '

# Row 25
$ws.Range('E25').Value = '  +1.58%  '

# Row 26
$ws.Range('D26').Value = '171.08'
$ws.Range('E26').Value = '  +4.45%  '

# Row 27
$ws.Range('D27').Value = '9.27'
$ws.Range('E27').Value = '  +8.26%  '

# Row 28
$ws.Range('D28').Value = '20.69'
$ws.Range('E28').Value = '  +5.54%  '

# Row 29
$ws.Range('E29').Value = '  +4.41%  '

# Row 30
$ws.Range('E30').Value = '  +1.11%  '

# Row 31
$ws.Range('E31').Value = '  +27.45%  '

# Row 32
$ws.Range('B32').Value = 'Gas'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D32').Value = '22.29'
$ws.Range('E32').Value = '  -1.09%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '4.52'
$ws.Range('E33').Value = '  +2.80%  '

# Row 34
$ws.Range('D34').Value = '0.0613'
$ws.Range('E34').Value = '  +4.77%  '

# Row 35
$ws.Range('D35').Value = '0.0908'
$ws.Range('E35').Value = '  +12.30%  '

# Row 36
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.17%  '

# Row 37
$ws.Range('D37').Value = '2.30'
$ws.Range('E37').Value = '  +6.13%  '

# Row 38
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '4.12'
$ws.Range('E38').Value = '  +3.37%  '

# Row 39
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').Value = '1.84'
$ws.Range('E39').Value = '  -1.28%  '

# Row 40
$ws.Range('E40').Value = '  +1.03%  '

# Row 41
$ws.Range('D41').Value = '18.33'
$ws.Range('E41').Value = '  +14.60%  '

# Row 42
$ws.Range('D42').Value = '0.0225'
$ws.Range('E42').Value = '  +4.14%  '

# Row 43
$ws.Range('E43').Value = '  +5.08%  '

# Row 44
$ws.Range('D44').Value = '98.83'
$ws.Range('E44').Value = '  +3.03%  '

# Row 45
$ws.Range('D45').Value = '0.0917'
$ws.Range('E45').Value = '  +12.04%  '

# Row 46
$ws.Range('E46').Value = '  +0.15%  '

# Row 47
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = '4.03'
$ws.Range('E47').Value = '  +98.03%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.319.52'
$ws.Range('E48').Value = '  +0.40%  '

# Row 49
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = '2.96'
$ws.Range('E49').Value = '  +7.29%  '

# Row 50
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = '7.05'
$ws.Range('E50').Value = '  +14.65%  '

# Row 51
$ws.Range('D51').Value = '2.292.32'
$ws.Range('E51').Value = '  +4.20%  '
